$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.358.25"
$ws.Range("E2").Value = "  -0.39%  "

$ws.Range("D3").Value = "2.625.66"
$ws.Range("E3").Value = "  -2.06%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "'594.96"
$ws.Range("E5").Value = "  -0.89%  "

$ws.Range("D6").Value = "'166.83"
$ws.Range("E6").Value = "  +0.52%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("E8").Value = "  -2.48%  "

$ws.Range("D9").Value = "2.625.34"
$ws.Range("E9").Value = "  -2.06%  "

$ws.Range("E10").Value = "  -2.41%  "

$ws.Range("E11").Value = "  +1.22%  "

$ws.Range("E12").Value = "  +1.47%  "

$ws.Range("D14").Value = "'27.66"
$ws.Range("E14").Value = "  -0.69%  "

$ws.Range("D15").Value = "3.103.25"
$ws.Range("E15").Value = "  -2.06%  "

$ws.Range("D16").Value = "'0.0000183"
$ws.Range("E16").Value = "  -0.97%  "

$ws.Range("D17").Value = "67.117.03"
$ws.Range("E17").Value = "  -0.64%  "

$ws.Range("D18").Value = "2.620.29"
$ws.Range("E18").Value = "  -1.92%  "

$ws.Range("E19").Value = "  +2.40%  "

$ws.Range("D20").Value = "'7.97"
$ws.Range("E20").Value = "  +4.20%  "

$ws.Range("D21").Value = "'357.67"
$ws.Range("E21").Value = "  -1.83%  "

$ws.Range("E22").Value = "  -1.50%  "

$ws.Range("E23").Value = "  -3.25%  "

$ws.Range("E24").Value = "  +0.06%  "

$ws.Range("E25").Value = "  -5.16%  "

$ws.Range("E26").Value = "  +1.32%  "

$ws.Range("D27").Value = "'69.80"
$ws.Range("E27").Value = "  -1.97%  "

$ws.Range("E29").Value = "  +0.02%  "

$ws.Range("E30").Value = "  -2.11%  "

$ws.Range("D31").Value = "'546.26"
$ws.Range("E31").Value = "  -2.16%  "

$ws.Range("D32").Value = "'7.93"
$ws.Range("E32").Value = "  -1.09%  "

$ws.Range("E33").Value = "  -3.16%  "

$ws.Range("E34").Value = "  -2.23%  "

$ws.Range("E35").Value = "  +4.45%  "

$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  +0.12%  "

$ws.Range("E37").Value = "  -3.56%  "

$ws.Range("D38").Value = "'156.70"
$ws.Range("E38").Value = "  +0.20%  "

$ws.Range("D39").Value = "'19.01"
$ws.Range("E39").Value = "  -2.92%  "

$ws.Range("E40").Value = "  -2.18%  "

$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'1.81"
$ws.Range("E41").Value = "  -1.19%  "

$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D42").Value = "'5.22"
$ws.Range("E42").Value = "  -1.84%  "

$ws.Range("E43").Value = "  +1.28%  "

$ws.Range("E44").Value = "  +0.09%  "

$ws.Range("E45").Value = "  -4.64%  "

$ws.Range("E46").Value = "  -0.34%  "

$ws.Range("D47").Value = "'152.35"
$ws.Range("E47").Value = "  -0.77%  "

$ws.Range("E48").Value = "  -1.92%  "

$ws.Range("D49").Value = "'3.78"
$ws.Range("E49").Value = "  -1.45%  "

$ws.Range("D50").Value = "'1.70"
$ws.Range("E50").Value = "  -1.62%  "

$ws.Range("D51").Value = "'0.0771"
$ws.Range("E51").Value = "  -0.87%  "
